$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '28.643.46'
Set-TextValue 'E2' '  +0.98%  '
Set-TextValue 'D3' '1.802.00'
Set-TextValue 'E3' '  -0.81%  '
Set-TextValue 'D4' '1.002'
Set-TextValue 'E4' '  +0.26%  '
Set-TextValue 'D5' '316.31'
Set-TextValue 'E5' '  -0.54%  '
Set-TextValue 'E6' '  +0.23%  '
Set-TextValue 'D7' '0.5305'
Set-TextValue 'E7' '  -8.56%  '
Set-TextValue 'D8' '0.3770'
Set-TextValue 'E8' '  -2.59%  '
Set-TextValue 'D9' '42.50'
Set-TextValue 'E9' '  -1.90%  '
Set-TextValue 'D10' '0.07499'
Set-TextValue 'E10' '  -1.67%  '
Set-TextValue 'D11' '1.114'
Set-TextValue 'E11' '  -2.29%  '
Set-TextValue 'D12' '1.002'
Set-TextValue 'E12' '  +0.27%  '
Set-TextValue 'D13' '20.71'
Set-TextValue 'E13' '  -2.61%  '
Set-TextValue 'D14' '6.149'
Set-TextValue 'E14' '  -1.86%  '
Set-TextValue 'D15' '7.342'
Set-TextValue 'E15' '  +0.49%  '
Set-TextValue 'D16' '1.798.40'
Set-TextValue 'E16' '  -0.98%  '
Set-TextValue 'D17' '90.30'
Set-TextValue 'E17' '  -2.13%  '
Set-TextValue 'E18' '  -1.39%  '
Set-TextValue 'D19' '0.06471'
Set-TextValue 'E19' '  -0.62%  '
Set-TextValue 'E20' '  +0.23%  '
Set-TextValue 'D21' '17.23'
Set-TextValue 'E21' '  -0.58%  '
Set-TextValue 'D22' '5.897'
Set-TextValue 'E22' '  -1.55%  '
Set-TextValue 'D23' '28.639.27'
Set-TextValue 'E23' '  +0.83%  '
Set-TextValue 'D24' '11.10'
Set-TextValue 'E24' '  -2.25%  '
Set-TextValue 'D25' '2.095'
Set-TextValue 'E25' '  -0.22%  '
Set-TextValue 'D26' '159.60'
Set-TextValue 'E26' '  +1.25%  '
Set-TextValue 'D27' '20.43'
Set-TextValue 'E27' '  -2.14%  '
Set-TextValue 'D28' '1.999.73'
Set-TextValue 'E28' '  -1.27%  '
Set-TextValue 'D29' '2.344'
Set-TextValue 'E29' '  -2.89%  '
Set-TextValue 'D30' '122.59'
Set-TextValue 'E30' '  -1.01%  '
Set-TextValue 'D31' '1.102'
Set-TextValue 'E31' '  -5.35%  '
Set-TextValue 'D32' '0.1053'
Set-TextValue 'E32' '  -1.32%  '
Set-TextValue 'B33' 'Filecoin'
Set-TextValue 'C33' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D33' '5.638'
Set-TextValue 'E33' '  -2.50%  '
Set-TextValue 'B34' 'HuobiToken'
Set-TextValue 'C34' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D34' '3.691'
Set-TextValue 'E34' '  +1.74%  '
Set-TextValue 'D35' '0.2240'
Set-TextValue 'E35' '  +3.41%  '
Set-TextValue 'E36' '  +5.16%  '
Set-TextValue 'D37' '0.02307'
Set-TextValue 'E37' '  -0.42%  '
Set-TextValue 'D38' '8.793'
Set-TextValue 'E38' '  -0.92%  '
Set-TextValue 'D39' '5.039'
Set-TextValue 'E39' '  -0.18%  '
Set-TextValue 'D40' '1.207'
Set-TextValue 'E40' '  +4.45%  '
Set-TextValue 'D41' '11.24'
Set-TextValue 'E41' '  -4.02%  '
Set-TextValue 'D42' '0.6205'
Set-TextValue 'E42' '  -3.40%  '
Set-TextValue 'E43' '  +0.29%  '
Set-TextValue 'D44' '1.411'
Set-TextValue 'E44' '  +2.08%  '
Set-TextValue 'D45' '13.23'
Set-TextValue 'E45' '  -1.27%  '
Set-TextValue 'D46' '3.690'
Set-TextValue 'E46' '  -0.54%  '
Set-TextValue 'D47' '0.5844'
Set-TextValue 'E47' '  -2.35%  '
Set-TextValue 'D48' '125.78'
Set-TextValue 'E48' '  +2.92%  '
Set-TextValue 'D49' '1.938'
Set-TextValue 'E49' '  -0.27%  '
Set-TextValue 'D50' '1.149'
Set-TextValue 'E50' '  +0.12%  '
Set-TextValue 'D51' '0.06889'
Set-TextValue 'E51' '  +0.57%  '
